# Add files via upload
# Added a checker for name/difficulty.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update row 2 (Ball brick collision) with notes about the brick-strength fix ---
$ws.Range("C2").Value = "Would bounce of both of and break or weaken them in oppisote direcion. BrickStrenth -1;"
$ws.Range("D2").Value = "When the ball hits two bircks it breakes them but continues going forward. Brick deletion and strenth changing works. Eg Brick strengh from 1 to 0"

# --- Row 4 Fix column picked up a trailing space ---
$ws.Range("C4").Value = "if they are down to 1 strenght when brick hits they will get deleted out of the array "

# --- Fill in the previously blank row 9 with the new "Enter Name" testing log entry ---
$ws.Range("A9").Value = "Enter Name"
$ws.Range("B9").Value = "Typing and enterbutton/clicking ok"
$ws.Range("C9").Value = "Will keep asking name until user types a name that is more than 1 character and not numbers"
$ws.Range("D9").Value = "Works when user clicks cancel or null or 1 character but not when they type numbers."
$ws.Range("E9").Value = 'Added !isNan(name) to the while loop that I already made check if name is = to "null" or empty.'
$ws.Rows("9").RowHeight = 33

# --- Update the selected cell to match the saved selection ---
$ws.Range("D2").Select()
